$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (slash format -> dash format).
# Leading apostrophe forces Excel to treat the value as literal text
# instead of auto-parsing ambiguous day<=12 strings as dates; resetting
# the style back to "Normal" afterwards clears the quote-prefix / number
# format residue so the cell ends up as a plain text value, matching the
# original (unstyled) inline-string cells.
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    $cell.Value = "'" + $dates[$row]
    $cell.Style = "Normal"
}

# Update the attendance counts for rows 3 and 5: D (Total Attendance Count)
# and G (Invalid) change from 0 to 1.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("G5").Value = 1
